# Applies the "vsc_trunk_name" / segmentation id additions to the "Vscs" sheet:
#   1. Inserts a new row 54 "VSC Trunk Name" (pushing the OpenStack external
#      network/subnet/port rows etc. down by one).
#   2. Inserts three new rows for "First/Second/Third External Port's
#      Segmentation Id" just before "OpenStack Port Name" (after the shift,
#      this lands at rows 75-77).
#   3. Adds comments describing the new fields, matching the style of the
#      existing sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vscs")

# --- Step 1: insert "VSC Trunk Name" row at row 54 -------------------------
$ws.Rows.Item(54).Insert()
# Row insert clones the style of the row it split from but puts the same
# style on every column; restore B/C formatting (border style, s="7") by
# copying formats from the row directly below (which still has the original
# "OpenStack first External Network" row's look).
$ws.Range("B55:C55").Copy()
$ws.Range("B54:C54").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A54").Value = "VSC Trunk Name"
$cmt1 = $ws.Range("A54").AddComment("Name of the trunk to be setup between control port and its underlay ports")

# --- Step 2: insert the three segmentation id rows before row 75 -----------
# After step 1, "OpenStack Port Name" / "OpenStack Port Security Groups" sit
# at rows 74/75 and "Third External Netmask Prefix Length" sits at row 74... 
# actually Port Name is 74 and Port Security Groups is 75; the new rows need
# to land above "OpenStack Port Security Groups" (row 75) so that, after the
# insert, the order becomes:
#   74 Third External Netmask Prefix Length (unchanged)
#   75 First External Port's Segmentation Id   (new)
#   76 Second External Port's Segmentation Id  (new)
#   77 Third External Port's Segmentation Id   (new)
#   78 OpenStack Port Name                      (shifted)
#   79 OpenStack Port Security Groups           (shifted)
$ws.Range("A75:A77").EntireRow.Insert()
$ws.Range("B78:C78").Copy()
$ws.Range("B75:C77").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A75").Value = "First External Port's Segmentation Id"
$ws.Range("A76").Value = "Second External Port's Segmentation Id"
$ws.Range("A77").Value = "Third External Port's Segmentation Id"

$cmt2 = $ws.Range("A75").AddComment("Segmentation id of first external port to be used during trunking")
$cmt3 = $ws.Range("A76").AddComment("Segmentation id of first external port to be used during trunking")
$cmt4 = $ws.Range("A77").AddComment("Segmentation id of third external port to be used during trunking")

# Integer data validation (matching the style used elsewhere on this sheet)
# for the three new "segmentation id" value cells (columns B and C).
foreach ($addr in @("B75","C75","B76","C76","B77","C77")) {
    $rng = $ws.Range($addr)
    $rng.Validation.Delete()
    $rng.Validation.Add(1, 2, 1, "")
    $rng.Validation.ErrorTitle = "Invalid Entry"
    $rng.Validation.ErrorMessage = "Your entry is not an integer, change anyway?"
    $rng.Validation.InputTitle = "Integer Selection"
    $rng.Validation.InputMessage = "Please provide integer"
    $rng.Validation.ShowError = $true
    $rng.Validation.ShowInput = $true
    $rng.Validation.IgnoreBlank = $true
}
